# Cotações atualizadas - 2025-12-01
# Append a new row (87) with the latest quotation data, mirroring the
# existing rows' layout: col A = numeric date serial, cols B:E = text
# values (comma-decimal, Portuguese locale) stored as plain strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

$ws.Cells.Item($row, 1).Value = 45992
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 2).Value = "21,7481"
$ws.Cells.Item($row, 3).Value = "10,9655"
$ws.Cells.Item($row, 4).Value = "15,5452"
$ws.Cells.Item($row, 5).Value = "15,5452"
